# Fix login selector issue:
#  - Update the Alphabets value to the new comma-separated list
#  - Flip Fresh_Run from Yes to No (bot already downloaded patient details)
#  - Fill in the File_Path with the location of the already-downloaded report

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

$ws.Range("B6").Value = "ge, gh, gi, gn, go, gr, gu, gw, ha"
$ws.Range("B7").Value = "No"
$ws.Range("B8").Value = "C:\Users\User\Desktop\Damco RPA Projects\Morgan_Records_Process_Wave_2\Output\Morgan Records Report 17 Nov 2020.xlsx"

$ws.Range("B9").Select()
